# Updated cryptos list values (Price & Volume(1h)) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed as a number by Excel;
# force them to Text format first so the literal string is preserved, matching the
# inline-string cell type used throughout column D.
$textPriceCells = @("D5", "D6", "D13", "D14", "D18", "D21", "D22", "D23", "D25", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D40", "D41", "D44", "D45", "D46", "D47", "D49", "D50")
foreach ($cellAddr in $textPriceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.295.44"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.782.24"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "594.32"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "167.54"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "3.783.50"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "36.11"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "4.416.42"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "3.772.56"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "68.292.72"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "17.90"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D21").Value = "10.77"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "464.86"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "0.698"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +8.41%  "
$ws.Range("D25").Value = "83.87"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "7.30"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "29.95"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("D34").Value = "9.14"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "3.736.47"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "43.77"
$ws.Range("E44").Value = "  +14.67%  "
$ws.Range("D45").Value = "0.300"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "47.08"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "146.19"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "392.67"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "2.786.00"
$ws.Range("E51").Value = "  +4.15%  "
